# "minor amendments in main_page set up"
#
# This script reworks the parameter table on the "main_page" sheet:
#  - swaps the "value"/"name" header columns and inserts a new "value"
#    sub-type column (value_box / mean_box / checkbox_group) ahead of the
#    existing "name" column (myvaluebox1 / myvaluebox2 / myvaluebox3 / ...)
#  - drops the now-unused value-box rows (myvaluebox4-8)
#  - adds a new "filter_box" / "checkbox_group" row (checkboxfilter1) that
#    drives a "User Enrolled" checkbox filter
#  - tidies up column widths / selection on main_page and demographics

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("main_page")
$ws3 = $wb.Worksheets.Item("demographics")

# ---- main_page : header row ------------------------------------------
$ws.Range("A1").Value = "type"
$ws.Range("B1").Value = "value"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "parameter_list"
$ws.Range("E1").Value = "variable"
$ws.Range("F1").Value = "variable_value"

# ---- row 2 : myvaluebox1 (joined) -------------------------------------
$ws.Range("A2").Value = "value_box"
$ws.Range("B2").Value = "value_box"
$ws.Range("C2").Value = "myvaluebox1"
$ws.Range("D2").Value = 'text = "Joined", colour = "yellow", icon = "user"'
$ws.Range("E2").Value = "joined"
$ws.Range("F2").Value = 1

# ---- row 3 : myvaluebox2 (enrolled) -----------------------------------
$ws.Range("A3").Value = "value_box"
$ws.Range("B3").Value = "value_box"
$ws.Range("C3").Value = "myvaluebox2"
$ws.Range("D3").Value = 'text = "Enrolled", colour = "purple", icon = "clipboard"'
$ws.Range("E3").Value = "enrolled"
$ws.Range("F3").ClearContents()

# ---- row 4 : myvaluebox3 -> mean_box (average days in chatbot) --------
$ws.Range("A4").Value = "value_box"
$ws.Range("B4").Value = "mean_box"
$ws.Range("C4").Value = "myvaluebox3"
$ws.Range("D4").Value = 'text = "Average days in the chatbot", colour = "orange", icon = "active"'
$ws.Range("E4").Value = "time_in_study_n"
$ws.Range("F4").ClearContents()

# ---- drop the old myvaluebox4-8 rows (now rows 5-9) --------------------
$ws.Rows("5:9").Delete()

# ---- new row 5 : filter_box / checkbox_group (checkboxfilter1) --------
$ws.Range("A5").Value = "filter_box"
$ws.Range("B5").Value = "checkbox_group"
$ws.Range("C5").Value = "checkboxfilter1"
$ws.Range("D5").Value = 'label = "User Enrolled", choices = c("yes", "no"), selected = c("yes", "no")'
$ws.Range("E5").Value = "enrolled"

# ---- column widths (approximate "best fit" values) ---------------------
$ws.Columns("A").ColumnWidth = 14.54296875 - (5/6)
$ws.Columns("B").ColumnWidth = 14.54296875 - (5/6)
$ws.Columns("C").ColumnWidth = 13.54296875 - (5/6)
$ws.Columns("D").ColumnWidth = 61.7265625  - (5/6)
$ws.Columns("E").ColumnWidth = 24.7265625  - (5/6)
$ws.Columns("F").ColumnWidth = 13          - (5/6)

# ---- page setup : portrait orientation ---------------------------------
$ws.PageSetup.Orientation = 1

# ---- demographics : new column width for column E ----------------------
$ws3.Columns("E").ColumnWidth = 43.90625 - (5/6)

# ---- restore selections (main_page stays the active sheet/tab) --------
$ws3.Activate()
$ws3.Range("E2").Select()

$ws.Activate()
$ws.Range("E16").Select()
